# Refresh the cryptos price list (Price / Volume(1h) columns) to match the
# latest scrape, and swap the Stacks / ApeXProtocol rows (47 <-> 48) to match
# the new ranking order, per the GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.057.07"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").Value = "2.495.20"
$ws.Range("E3").Value = "  -0.53%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'320.23"
$ws.Range("E5").Value = "  -0.77%  "

$ws.Range("D6").Value = "'107.33"
$ws.Range("E6").Value = "  -2.41%  "

$ws.Range("D7").Value = "'0.525"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "'0.536"
$ws.Range("E9").Value = "  -2.05%  "

$ws.Range("D10").Value = "'39.45"
$ws.Range("E10").Value = "  -1.72%  "

$ws.Range("D11").Value = "'20.15"
$ws.Range("E11").Value = "  +7.94%  "

$ws.Range("D12").Value = "'0.0810"
$ws.Range("E12").Value = "  -0.39%  "

$ws.Range("E13").Value = "  -0.16%  "

$ws.Range("D14").Value = "'7.11"
$ws.Range("E14").Value = "  -1.99%  "

$ws.Range("D15").Value = "2.891.47"
$ws.Range("E15").Value = "  -0.18%  "

$ws.Range("D16").Value = "2.503.62"
$ws.Range("E16").Value = "  +0.10%  "

$ws.Range("D17").Value = "'0.835"
$ws.Range("E17").Value = "  -2.05%  "

$ws.Range("D18").Value = "47.931.80"
$ws.Range("E18").Value = "  +0.12%  "

$ws.Range("D19").Value = "'12.94"
$ws.Range("E19").Value = "  -3.13%  "

$ws.Range("D20").Value = "'6.74"
$ws.Range("E20").Value = "  +1.15%  "

$ws.Range("D21").Value = "0.0₃0936"
$ws.Range("E21").Value = "  -1.12%  "

$ws.Range("D22").Value = "'2.72"
$ws.Range("E22").Value = "  -1.29%  "

$ws.Range("D23").Value = "'276.22"
$ws.Range("E23").Value = "  +11.38%  "

$ws.Range("D24").Value = "'71.71"
$ws.Range("E24").Value = "  +1.25%  "

$ws.Range("D25").Value = "'2.55"
$ws.Range("E25").Value = "  -0.45%  "

$ws.Range("D27").Value = "'25.60"
$ws.Range("E27").Value = "  -1.47%  "

$ws.Range("D28").Value = "'9.71"
$ws.Range("E28").Value = "  -3.37%  "

$ws.Range("E29").Value = "  +0.21%  "

$ws.Range("D30").Value = "'34.93"
$ws.Range("E30").Value = "  -0.31%  "

$ws.Range("D31").Value = "'2.10"
$ws.Range("E31").Value = "  -4.59%  "

$ws.Range("D32").Value = "'49.38"
$ws.Range("E32").Value = "  -1.18%  "

$ws.Range("D33").Value = "'19.57"
$ws.Range("E33").Value = "  -3.47%  "

$ws.Range("D34").Value = "'1.01"
$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("D35").Value = "'5.29"
$ws.Range("E35").Value = "  -1.40%  "

$ws.Range("D36").Value = "'0.0777"
$ws.Range("E36").Value = "  -1.70%  "

$ws.Range("D37").Value = "'1.94"
$ws.Range("E37").Value = "  -2.08%  "

$ws.Range("D38").Value = "'4.61"
$ws.Range("E38").Value = "  -2.28%  "

$ws.Range("D39").Value = "'2.88"
$ws.Range("E39").Value = "  -2.58%  "

$ws.Range("E40").Value = "  -0.93%  "

$ws.Range("D41").Value = "'121.10"
$ws.Range("E41").Value = "  +1.16%  "

$ws.Range("E42").Value = "  -0.26%  "

$ws.Range("D43").Value = "'21.33"
$ws.Range("E43").Value = "  -5.47%  "

$ws.Range("D44").Value = "'0.0300"
$ws.Range("E44").Value = "  +0.07%  "

$ws.Range("D45").Value = "2.008.80"
$ws.Range("E45").Value = "  +0.45%  "

$ws.Range("D46").Value = "'3.14"
$ws.Range("E46").Value = "  +2.51%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'2.00"
$ws.Range("E47").Value = "  -2.27%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'1.84"
$ws.Range("E48").Value = "  -0.80%  "

$ws.Range("D49").Value = "'8.98"
$ws.Range("E49").Value = "  -1.06%  "

$ws.Range("D50").Value = "'5.15"
$ws.Range("E50").Value = "  -1.88%  "

$ws.Range("D51").Value = "'79.99"
$ws.Range("E51").Value = "  +2.40%  "
